$wb = $excel.ActiveWorkbook

# Rename the existing sheet "Sheet1" -> "1d"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "1d"

# Insert a new worksheet right after "1d" and name it "2d"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "2d"

# Populate the new "2d" sheet with the recomputed comparison results
$ws2.Range("B1").Value = "dice coeff"

$ws2.Range("A3").Value = "Monte-Carlo"
$ws2.Range("B3").Value = "IsoContour"
$ws2.Range("C3").Value = "ImgSum"

$ws2.Range("A4").Value = 100
$ws2.Range("B4").Value = 0.99091304242834
$ws2.Range("C4").Value = 0.990670084076037

$ws2.Range("A5").Value = 200
$ws2.Range("B5").Value = 0.993180682270621
$ws2.Range("C5").Value = 0.992987704768507

$ws2.Range("A6").Value = 400
$ws2.Range("B6").Value = 0.995458724609977
$ws2.Range("C6").Value = 0.992832473889901

$ws2.Range("A7").Value = 800
$ws2.Range("B7").Value = 0.995412844036697
$ws2.Range("C7").Value = 0.999580696433969

$ws2.Range("A10").Value = "Collocation"
$ws2.Range("B10").Value = "IsoContour"

$ws2.Range("A11").Value = 3
$ws2.Range("B11").Value = 0.995449065190443

$ws2.Range("A12").Value = 4
$ws2.Range("B12").Value = 0.995447857776182

$ws2.Range("A13").Value = 5
$ws2.Range("B13").Value = 0.995434576412625

# Leave the cursor where the author left it on the new sheet
[void]$ws2.Range("B14").Select()

# "2d" becomes the active/selected sheet (matches activeTab=1, tabSelected flags)
$ws2.Activate()
